$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Croatia ✓ - Gibraltar: 3:0"
$ws.Range("G2").Value = "✓"

$ws.Range("A3").Value = "Netherlands ✓ - Finland: 4:0"
$ws.Range("G3").Value = "✓"

$ws.Range("A4").Value = "Scotland ✓ - Belarus: 2:1"
$ws.Range("G4").Value = "✓"

$ws.Range("A5").Value = "Romania - Austria X: 1:0"
$ws.Range("G5").Value = "X"

$ws.Range("A6").Value = "Egypt ✓ - Guinea-Bissau: 1:0"
$ws.Range("G6").Value = "✓"

$ws.Range("A7").Value = "Burkina Faso ✓ - Ethiopia: 3:1"
$ws.Range("G7").Value = "✓"

$ws.Range("A8").Value = "Ghana ✓ - Comoros: 1:0"
$ws.Range("G8").Value = "✓"

$ws.Range("A9").Value = "Club Deportivo Guabirá ✓ - Club Aurora: 2:1"
$ws.Range("G9").Value = "✓"
